$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.469.73'
$ws.Range("E2").Value = '  +0.99%  '
$ws.Range("D3").Value = '1.879.21'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.76'
$ws.Range("E5").Value = '  +5.29%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4761'
$ws.Range("E7").Value = '  +1.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2901'
$ws.Range("E8").Value = '  +1.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06532'
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.84'
$ws.Range("E10").Value = '  +3.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07743'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.89'
$ws.Range("E12").Value = '  +3.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7368'
$ws.Range("E13").Value = '  +7.61%  '
$ws.Range("D14").Value = '1.879.78'
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.130'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '273.09'
$ws.Range("E16").Value = '  +1.32%  '
$ws.Range("D17").Value = '30.462.58'
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.64'
$ws.Range("E18").Value = '  +2.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007604'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").Value = '2.125.72'
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.244'
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.183'
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.306'
$ws.Range("E25").Value = '  -0.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.20'
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.91'
$ws.Range("E28").Value = '  +2.49%  '
$ws.Range("E29").Value = '  +0.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09968'
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.518'
$ws.Range("E31").Value = '  +4.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.321'
$ws.Range("E32").Value = '  +2.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.070'
$ws.Range("E33").Value = '  +1.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04777'
$ws.Range("E34").Value = '  +1.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.126'
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7009'
$ws.Range("E36").Value = '  +1.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.715'
$ws.Range("E37").Value = '  +0.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01869'
$ws.Range("E38").Value = '  +1.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.731'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.339'
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.939'
$ws.Range("E41").Value = '  +2.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.40'
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4195'
$ws.Range("E43").Value = '  +3.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.0000'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8398'
$ws.Range("E45").Value = '  +0.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.75'
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.255'
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.088'
$ws.Range("E48").Value = '  +1.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.53'
$ws.Range("E49").Value = '  +4.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '927.93'
$ws.Range("E50").Value = '  -1.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05639'
$ws.Range("E51").Value = '  +1.13%  '
